# Refresh the scraped cryptocurrency price/volume snapshot (GitHub Actions run).
#
# The source data are plain text cells (prices can contain multiple "."
# thousands separators, e.g. "42.424.64", so they are never real numbers).
# Some of the new Price values, though, DO look like plain decimals (e.g.
# "300.40", "1.00"), and Excel's normal Range.Value auto-detection would
# silently coerce those into numeric cells. To keep them as literal text we
# enter them with a leading apostrophe (the standard "force text" quote
# prefix) and then reset the cell style back to Normal so no stray
# quote-prefix formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = '42.424.64'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '2.291.25'
$ws.Range("E3").Value = '  +0.30%  '
Set-TextValue "D4" '1.00'
$ws.Range("E4").Value = '  -0.02%  '
Set-TextValue "D5" '300.40'
$ws.Range("E5").Value = '  -1.52%  '
Set-TextValue "D6" '94.80'
$ws.Range("E6").Value = '  -0.70%  '
Set-TextValue "D7" '0.506'
$ws.Range("E7").Value = '  +0.64%  '
$ws.Range("E8").Value = '  -0.01%  '
Set-TextValue "D9" '0.490'
$ws.Range("E9").Value = '  -1.41%  '
Set-TextValue "D10" '34.33'
$ws.Range("E10").Value = '  -1.95%  '
Set-TextValue "D11" '18.91'
$ws.Range("E11").Value = '  +3.63%  '
Set-TextValue "D12" '0.0779'
$ws.Range("E12").Value = '  -0.58%  '
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").Value = '2.646.80'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").Value = '2.303.34'
$ws.Range("E16").Value = '  +0.29%  '
Set-TextValue "D17" '0.776'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").Value = '42.374.76'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("E19").Value = '  -5.55%  '
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("E21").Value = '  -0.45%  '
Set-TextValue "D22" '67.46'
$ws.Range("E22").Value = '  +0.94%  '
Set-TextValue "D23" '235.06'
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("E24").Value = '  +5.87%  '
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("E26").Value = '  -1.98%  '
$ws.Range("E27").Value = '  -3.10%  '
$ws.Range("E28").Value = '  +4.45%  '
Set-TextValue "D29" '164.30'
$ws.Range("E29").Value = '  -0.76%  '
Set-TextValue "D30" '9.00'
$ws.Range("E30").Value = '  +0.29%  '
Set-TextValue "D31" '31.54'
$ws.Range("E31").Value = '  -2.93%  '
Set-TextValue "D32" '1.00'
$ws.Range("E32").Value = '  -0.05%  '
Set-TextValue "D33" '4.96'
$ws.Range("E33").Value = '  +0.60%  '
Set-TextValue "D34" '17.36'
$ws.Range("E34").Value = '  -0.28%  '
$ws.Range("E35").Value = '  +0.54%  '
$ws.Range("E36").Value = '  -2.51%  '
$ws.Range("E37").Value = '  -7.59%  '
Set-TextValue "D38" '0.0993'
$ws.Range("E38").Value = '  -1.29%  '
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("E41").Value = '  -0.04%  '
Set-TextValue "D42" '19.67'
$ws.Range("E42").Value = '  +9.56%  '
$ws.Range("D43").Value = '1.943.89'
$ws.Range("E43").Value = '  -2.19%  '
Set-TextValue "D44" '10.28'
$ws.Range("E44").Value = '  +3.10%  '
$ws.Range("E45").Value = '  -0.21%  '
$ws.Range("E46").Value = '  +3.58%  '
Set-TextValue "D47" '2.72'
$ws.Range("E47").Value = '  -1.68%  '
$ws.Range("B48").Value = 'HuobiToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D48" '2.83'
$ws.Range("E48").Value = '  -1.61%  '
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.515.53'
$ws.Range("E49").Value = '  +0.26%  '
Set-TextValue "D50" '52.66'
$ws.Range("E50").Value = '  -1.36%  '
$ws.Range("E51").Value = '  +0.81%  '
